# "Missed sound, fixed it"
# Insert a new row for the missing "SoundAmbientBirds" sound entry in the
# Scene section's Sounds list (after SoundMainTheme / SoundGameMusic), which
# pushes all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Knight+Player")

# Insert a new row above row 33 (SoundGameMusic is row 31, row 32 was blank,
# "Scripts" header was row 33) - shifting everything from row 33 downward by one.
$ws.Rows("33:33").Insert()

# Fill in the newly created row with the missing sound asset name.
$ws.Cells.Item(32, 2).Value = "SoundAmbientBirds"

# Restore the view/selection state recorded for this sheet after the edit.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("E25").Select()
